# Applies the cell-content edits described by the commit diff
# ("Updated symbol list on Mon Jan  9 07:37:51 UTC 2023 with GitHub Actions")
# to the cryptocurrency table on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text, and whether the text must be
# forced to stay text (numeric-looking values like '278.31' or '6.59%'
# would otherwise be auto-converted to numbers by Excel, which would
# silently drop significant trailing zeros, e.g. "0.8800" -> 0.88).
$edits = @(
    [PSCustomObject]@{ Cell = 'D2'; Value = '278.31'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E2'; Value = '6.59%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D3'; Value = '27.41'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E3'; Value = '1.47%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D4'; Value = '4.808'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E4'; Value = '1.91%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D5'; Value = '0.06238'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E5'; Value = '0.31%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D6'; Value = '6.921'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E6'; Value = '2.59%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'B7'; Value = 'MXToken'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'C7'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'D7'; Value = '0.8800'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E7'; Value = '3.46%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'B8'; Value = 'FTXToken'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'C8'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'D8'; Value = '0.9426'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E8'; Value = '3.04%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'B9'; Value = 'WazirX'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'C9'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'D9'; Value = '0.1453'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E9'; Value = '3.52%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'B10'; Value = 'LiechtensteinCryptoassetsExchange'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'C10'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'D10'; Value = '0.05248'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E10'; Value = '3.40%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'B11'; Value = 'MandalaExchangeToken'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'C11'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'D11'; Value = '0.07331'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E11'; Value = '3.63%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'B12'; Value = 'BitrueCoin'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'C12'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'D12'; Value = '0.03104'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E12'; Value = '0.19%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'B13'; Value = 'BitMartToken'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'C13'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'D13'; Value = '0.09058'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E13'; Value = '0.03%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'B14'; Value = 'BitForexToken'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'C14'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'D14'; Value = '0.001551'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E14'; Value = '1.47%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'B15'; Value = 'One'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'C15'; Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'D15'; Value = '0.0006266'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E15'; Value = '1.45%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'B16'; Value = 'TigerCash'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'C16'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'D16'; Value = '0.006016'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E16'; Value = '0.59%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'B17'; Value = 'LEO'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'C17'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'D17'; Value = '3.451'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E17'; Value = '0.21%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'B18'; Value = 'GateToken'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'C18'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; ForceText = $false }
    [PSCustomObject]@{ Cell = 'D18'; Value = '3.269'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E18'; Value = '3.17%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E19'; Value = '6.52%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D20'; Value = '0.3147'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E20'; Value = '1.30%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D21'; Value = '0.1303'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E21'; Value = '-0.60%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D22'; Value = '3.850'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E22'; Value = '-5.70%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D23'; Value = '0.04324'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E23'; Value = '1.74%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E24'; Value = '-2.00%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D25'; Value = '0.004281'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E25'; Value = '5.05%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D26'; Value = '0.0001200'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E26'; Value = '-0.03%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D27'; Value = '0.0001691'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E27'; Value = '3.11%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D40'; Value = '0.04044'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E40'; Value = '2.32%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D41'; Value = '0.006701'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E41'; Value = '62.20%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D42'; Value = '0.1155'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E42'; Value = '3.81%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D43'; Value = '0.002134'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E43'; Value = '-3.47%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D44'; Value = '0.01212'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E44'; Value = '-10.87%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D45'; Value = '0.00005107'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E45'; Value = '-1.06%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E46'; Value = '0.01%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E47'; Value = '852.86%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D49'; Value = '0.00002101'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E49'; Value = '0.01%'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'D50'; Value = '0.0002001'; ForceText = $true }
    [PSCustomObject]@{ Cell = 'E50'; Value = '0.01%'; ForceText = $true }
)

foreach ($edit in $edits) {
    $c = $ws.Range($edit.Cell)
    if ($edit.ForceText) {
        # Mark the cell as Text-formatted so Excel stores the literal
        # characters instead of parsing/re-formatting them as a number.
        $c.NumberFormat = "@"
    }
    $c.Value = $edit.Value
}
